# Apply edit: insert new "source_file" column at A, shift existing data
# columns right by one, update source_file values to new naming scheme,
# and refresh several rows with rerun content (per commit: "rerun to get new name").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = 'source_file'
$ws.Cells.Item(1, 2).Value = 'Latitude'
$ws.Cells.Item(1, 3).Value = 'Longitude'
$ws.Cells.Item(1, 4).Value = 'Date_of_Issuance'
$ws.Cells.Item(1, 5).Value = 'Condition_summary_1'
$ws.Cells.Item(1, 6).Value = 'Condition_summary_2'
$ws.Cells.Item(1, 7).Value = 'Condition_summary_3'
$ws.Cells.Item(1, 8).Value = 'Condition_summary_4'
$ws.Cells.Item(1, 9).Value = 'Habitat_Type'
$ws.Cells.Item(1, 10).Value = 'Fish_species'
$ws.Cells.Item(1, 11).Value = 'Offset_footprint_size'
$ws.Cells.Item(1, 12).Value = 'Vegetation_Cover'
$ws.Cells.Item(1, 13).Value = 'Boulder'
$ws.Cells.Item(1, 14).Value = 'Woody_coverage'
$ws.Cells.Item(1, 15).Value = 'Instream_structures'
$ws.Cells.Item(1, 16).Value = 'Langitude'

# Row 2
$ws.Cells.Item(2, 1).Value = 'OCR_18-HCAA-00233.json'
$ws.Cells.Item(2, 2).NumberFormat = "@"
$ws.Cells.Item(2, 2).Value = '44.1913'
$ws.Cells.Item(2, 3).NumberFormat = "@"
$ws.Cells.Item(2, 3).Value = '-78.8722'
$ws.Cells.Item(2, 4).Value = 'JAN 10 2020'
$ws.Cells.Item(2, 5).Value = 'The Proponent must install vegetation-enhanced armour stone walls on a 25 m section of the east bank (50 m²) and retrofit 10 parking lot catchbasins with shields along a 565 m section of the Moira River east bank.'
$ws.Cells.Item(2, 6).Value = 'Monitoring of offsetting measures for two years post-construction, with annual reports by November 30.'
$ws.Cells.Item(2, 7).Value = 'Prohibition on activities adversely impacting offsetting measures and requirement for access permissions.'
$ws.Cells.Item(2, 8).Value = 'Authorization limitations, including compliance with other regulations and prohibition on deleterious substance deposits.'
$ws.Cells.Item(2, 9).Value = 'river bank, parking lot catchbasins'
$ws.Cells.Item(2, 10).Value = 'Channel Darter'
$ws.Cells.Item(2, 11).Value = '50 m² vegetation-enhanced armour stone walls, 10 catchbasin shields over 565 m river section'
$ws.Cells.Item(2, 12).Value = '80% vegetation survival required for effectiveness'
$ws.Cells.Item(2, 13).Value = 'armour stone used in revetment'
$ws.Cells.Item(2, 14).Value = 'None explicitly mentioned'
$ws.Cells.Item(2, 15).Value = 'catchbasin shields'
$ws.Cells.Item(2, 16).Value = $null

# Row 3
$ws.Cells.Item(3, 1).Value = 'OCR_14-HCAA-00814.json'
$ws.Cells.Item(3, 2).Value = $null
$ws.Cells.Item(3, 3).Value = 'None'
$ws.Cells.Item(3, 4).Value = 'JUN 05 2015'
$ws.Cells.Item(3, 5).Value = 'development of an operating plan for the Timiskaming Dam to protect fish spawning habitats'
$ws.Cells.Item(3, 6).Value = 'monitoring and reporting requirements including annual reports on dam operations, fish recruitment, and contingency measures'
$ws.Cells.Item(3, 7).Value = 'implementation of adaptive management strategies if deviations from the operating plan impact spawning success'
$ws.Cells.Item(3, 8).Value = 'compliance with Species at Risk Act and prohibition on depositing deleterious substances'
$ws.Cells.Item(3, 9).Value = 'spawning habitats for Lake Sturgeon and Lake Whitefish'
$ws.Cells.Item(3, 10).Value = 'Lake Sturgeon, Lake Whitefish'
$ws.Cells.Item(3, 11).Value = 'None'
$ws.Cells.Item(3, 12).Value = 'None'
$ws.Cells.Item(3, 13).Value = 'None'
$ws.Cells.Item(3, 14).Value = 'None'
$ws.Cells.Item(3, 15).Value = 'None'
$ws.Cells.Item(3, 16).Value = 'None'

# Row 4
$ws.Cells.Item(4, 1).Value = 'OCR_14-HCAA-00810.json'
$ws.Cells.Item(4, 2).Value = 'None'
$ws.Cells.Item(4, 3).Value = 'None'
$ws.Cells.Item(4, 4).Value = 'may 08 2015'
$ws.Cells.Item(4, 5).Value = 'mitigation measures including fish rescue reporting, habitat offsetting, and compliance with SARA'
$ws.Cells.Item(4, 6).Value = 'monitoring and reporting requirements for both mitigation and offsetting measures, including annual reports with data and photographs'
$ws.Cells.Item(4, 7).Value = 'implementation of offsetting measures during construction phase with specific habitat enhancements'
$ws.Cells.Item(4, 8).Value = 'prohibition on transferring authorization and requirement to maintain on-site documentation'
$ws.Cells.Item(4, 9).Value = 'river, lake, shoreline, spawning, resting, foraging'
$ws.Cells.Item(4, 10).Value = 'walleye'
$ws.Cells.Item(4, 11).Value = '880 m2 (river habitat from old dam removal), 400 m2 (Walleye spawning and resting), 895 m2 (shoreline/upperland), 720 m2 (lake to river habitat)'
$ws.Cells.Item(4, 12).Value = 'None'
$ws.Cells.Item(4, 13).Value = 'additional boulders for Walleye spawning areas if washed out'
$ws.Cells.Item(4, 14).Value = 'None'
$ws.Cells.Item(4, 15).Value = 'rock shoal (5-10m length, 10-40 m2 area)'
$ws.Cells.Item(4, 16).Value = $null

# Row 5
$ws.Cells.Item(5, 1).Value = 'OCR_18-HCAA-00064.json'
$ws.Cells.Item(5, 2).Value = 'None'
$ws.Cells.Item(5, 3).Value = 'None'
$ws.Cells.Item(5, 4).Value = 'Feb 07/2020'
$ws.Cells.Item(5, 5).Value = 'The Proponent must implement offsetting measures including marsh habitat creation, tributary restoration, riparian buffer planting, and storm water management pond expansion to mitigate impacts on fish and fish habitat. Contingency measures are required if these measures fail to meet criteria.'
$ws.Cells.Item(5, 6).Value = 'Monitoring and reporting are mandatory, including post-construction assessments, fish sampling, and vegetation survival checks. Reports must be submitted by June 14, 2022.'
$ws.Cells.Item(5, 7).Value = 'Prohibition on depositing deleterious substances in water frequented by fish. Compliance with other regulatory agencies is required.'
$ws.Cells.Item(5, 8).Value = 'The Proponent is solely responsible for design and safety of works. Authorization cannot be transferred without prior notification to DFO.'
$ws.Cells.Item(5, 9).Value = 'marsh habitat, unnamed tributary restoration, riparian buffer, storm water management pond, realigned drain channels'
$ws.Cells.Item(5, 10).Value = 'None explicitly listed in section 4; general references to fish utilization and species composition in monitoring sections'
$ws.Cells.Item(5, 11).Value = 'marsh habitat expansion (size unspecified), unnamed tributary restoration (size unspecified), riparian buffer planting (80% coverage target), Hooper Drain channel (morphology retention), Central Drain channel (morphology retention), SWM pond (average depth 0.9m)'
$ws.Cells.Item(5, 12).Value = '80% survival for aquatic vegetation in marsh and 80% coverage for riparian plantings'
$ws.Cells.Item(5, 13).Value = 'None explicitly mentioned'
$ws.Cells.Item(5, 14).Value = 'Habitat enhancement features (e.g., woody debris) in marsh and tributary require no destabilization'
$ws.Cells.Item(5, 15).Value = 'Habitat enhancement features in marsh and tributary must remain stable'
$ws.Cells.Item(5, 16).Value = $null

# Row 6
$ws.Cells.Item(6, 1).Value = 'OCR_18-HCAA-00311.json'
$ws.Cells.Item(6, 2).Value = $null
$ws.Cells.Item(6, 3).Value = 'None'
$ws.Cells.Item(6, 4).Value = 'oct 16, 2019'
$ws.Cells.Item(6, 5).Value = 'authorization under Fisheries Act for work likely to cause serious harm to fish, requiring compliance with conditions including habitat offsetting, monitoring, and reporting'
$ws.Cells.Item(6, 6).Value = 'prohibition on depositing deleterious substances in water frequented by fish'
$ws.Cells.Item(6, 7).Value = 'requirement to obtain access permissions for lands/waters not owned by proponent'
$ws.Cells.Item(6, 8).Value = 'no transfer of authorization without prior notification to DFO'
$ws.Cells.Item(6, 9).Value = 'stream, channel, riparian zone, instream habitat features'
$ws.Cells.Item(6, 10).Value = 'rainbow trout'
$ws.Cells.Item(6, 11).Value = 'reach sn01, reach sn02, unnamed tributary of seneca creek (specific sizes not quantified in m2)'
$ws.Cells.Item(6, 12).Value = '80% survival target for planted vegetation in riparian zone'
$ws.Cells.Item(6, 13).Value = 'None'
$ws.Cells.Item(6, 14).Value = 'None'
$ws.Cells.Item(6, 15).Value = 'constructed instream habitat features including channel morphological features, spawning surveys required'
$ws.Cells.Item(6, 16).Value = 'None'

# Row 7
$ws.Cells.Item(7, 1).Value = 'OCR_18-HCAA-00192.json'
$ws.Cells.Item(7, 2).Value = '53°36''40.96"N'
$ws.Cells.Item(7, 3).Value = '108°44''38.01"W'
$ws.Cells.Item(7, 4).Value = 'FEB 26 2020'
$ws.Cells.Item(7, 5).Value = 'Sediment and erosion control measures must be in place and maintained to avoid sediment entering the NSR during all project phases.'
$ws.Cells.Item(7, 6).Value = 'Fish screens must meet requirements of 0.125 m/s intake flow with 7,045 m² of total screen area to prevent entrainment.'
$ws.Cells.Item(7, 7).Value = 'Installation of a gravel/cobble bar (5,800 m²) to offset habitat loss, using rounded stone 15–150 mm in diameter.'
$ws.Cells.Item(7, 8).Value = 'Monitoring and reporting requirements include post-construction surveys, embeddedness surveys, and functional monitoring over three years.'
$ws.Cells.Item(7, 9).Value = 'river bank, gravel/cobble bar'
$ws.Cells.Item(7, 10).Value = 'None explicitly listed in section 4 or elsewhere'
$ws.Cells.Item(7, 11).Value = '5800 m2 (gravel/boulder bar)'
$ws.Cells.Item(7, 12).Value = 'None specified'
$ws.Cells.Item(7, 13).Value = 'Gravel/boulder bar with stones 15–150 mm'
$ws.Cells.Item(7, 14).Value = 'None mentioned'
$ws.Cells.Item(7, 15).Value = 'Gravel/boulder bar installation'
$ws.Cells.Item(7, 16).Value = $null

# Row 8
$ws.Cells.Item(8, 1).Value = 'OCR_18-HCAA-00146.json'
$ws.Cells.Item(8, 2).Value = '5672412N'
$ws.Cells.Item(8, 3).Value = '11678490E'
$ws.Cells.Item(8, 4).Value = 'AUG 17 2018'
$ws.Cells.Item(8, 5).Value = 'Sedimentation and erosion control measures must be in place and maintained to avoid sediment release into the watercourse.'
$ws.Cells.Item(8, 6).Value = 'Total suspended sediment and turbidity monitoring must adhere to the plan in Appendix C of the 2018 Aquatic Effects Assessment.'
$ws.Cells.Item(8, 7).Value = 'All riprap must be clean, free of fine materials, and not obtained from fish-frequented waters below the ordinary high water mark.'
$ws.Cells.Item(8, 8).Value = 'Dewatering must direct water to vegetated areas or settling basins, ensuring water quality meets standards before returning to fish habitats.'
$ws.Cells.Item(8, 9).Value = 'fish habitat (specific type not explicitly stated beyond ''resident fish species'' and ''fish-frequented waters'')'
$ws.Cells.Item(8, 10).Value = 'resident fish species (specific species not listed)'
$ws.Cells.Item(8, 11).Value = '7,800 m² for the secondary channel regrading downstream'
$ws.Cells.Item(8, 12).Value = 'None explicitly mentioned beyond general requirements for sediment control and spoil disposal areas'
$ws.Cells.Item(8, 13).Value = 'Natural structures like large boulders contributing to fish habitat must be stockpiled for replacement post-construction'
$ws.Cells.Item(8, 14).Value = 'Natural structures including woody debris must be stockpiled for replacement'
$ws.Cells.Item(8, 15).Value = 'Berm construction and riprap placement are noted as in-stream activities with specific mitigation measures'
$ws.Cells.Item(8, 16).Value = $null

# Row 9
$ws.Cells.Item(9, 1).Value = 'OCR_18-HCAA-00145.json'
$ws.Cells.Item(9, 2).NumberFormat = "@"
$ws.Cells.Item(9, 2).Value = '43.79381'
$ws.Cells.Item(9, 3).NumberFormat = "@"
$ws.Cells.Item(9, 3).Value = '-80.386060'
$ws.Cells.Item(9, 4).Value = 'JAN 10 2020'
$ws.Cells.Item(9, 5).Value = 'Sediment and erosion control measures must be in place and maintained to avoid sediment release during work.'
$ws.Cells.Item(9, 6).Value = 'Fish rescue and relocation required before work in isolated areas.'
$ws.Cells.Item(9, 7).Value = 'Monitoring and reporting including photographic records and as-built surveys.'
$ws.Cells.Item(9, 8).Value = 'Habitat offsetting measures with specific criteria (e.g., overwintering pools, vegetation coverage).'
$ws.Cells.Item(9, 9).Value = 'fish habitat (permanent alteration, destruction)'
$ws.Cells.Item(9, 10).Value = 'trout'
$ws.Cells.Item(9, 11).Value = '50 m2 (boulder clusters), 100 m2 (sweeper trees), 80% vegetation coverage'
$ws.Cells.Item(9, 12).NumberFormat = "@"
$ws.Cells.Item(9, 12).Value = '80%'
$ws.Cells.Item(9, 13).Value = 'large boulder clusters (50 m2)'
$ws.Cells.Item(9, 14).Value = 'anchored sweeper trees (100 m2)'
$ws.Cells.Item(9, 15).Value = 'overwintering pool habitat (1.7m depth), boulder clusters, sweeper trees'
$ws.Cells.Item(9, 16).Value = $null

# Row 10
$ws.Cells.Item(10, 1).Value = 'OCR_18-HCAA-00253.json'
$ws.Cells.Item(10, 2).NumberFormat = "@"
$ws.Cells.Item(10, 2).Value = '50.894225'
$ws.Cells.Item(10, 3).NumberFormat = "@"
$ws.Cells.Item(10, 3).Value = '-114.009975'
$ws.Cells.Item(10, 4).Value = 'feb 04 2018'
$ws.Cells.Item(10, 5).Value = 'Sediment and erosion control measures must be in place and maintained to avoid sediment release into water. In-water activity timing restrictions to protect spawning fish and their eggs (May 1st to July 15th and September 16th to April 5th).'
$ws.Cells.Item(10, 6).Value = 'Berm construction and removal must adhere to approved plans, with footprint not exceeding 5% of design area without DFO approval. Natural structures must be replaced after construction.'
$ws.Cells.Item(10, 7).Value = 'Revegetation with native plants, equipment cleanliness to prevent invasive species, spill response plans, and stockpiling materials above high water level.'
$ws.Cells.Item(10, 8).Value = '3,462 m² of fish habitat restoration upstream/downstream on west bank. Monitoring and reporting on offset effectiveness over three years.'
$ws.Cells.Item(10, 9).Value = 'river habitat (offsetting measures on west bank)'
$ws.Cells.Item(10, 10).Value = 'None explicitly listed in the document'
$ws.Cells.Item(10, 11).Value = '3462 m²'
$ws.Cells.Item(10, 12).Value = 'None specified'
$ws.Cells.Item(10, 13).Value = 'None specified'
$ws.Cells.Item(10, 14).Value = 'None specified'
$ws.Cells.Item(10, 15).Value = 'None specified'
$ws.Cells.Item(10, 16).Value = $null

# Row 11
$ws.Cells.Item(11, 1).Value = 'OCR_18-HCAA-00160.json'
$ws.Cells.Item(11, 2).Value = '71.889403°N'
$ws.Cells.Item(11, 3).Value = '-80.887592°W'
$ws.Cells.Item(11, 4).Value = 'March 21, 2019'
$ws.Cells.Item(11, 5).Value = 'The work must be completed by the expiration date or DFO must be notified for extension.'
$ws.Cells.Item(11, 6).Value = 'Implement sediment and erosion control measures, including approved plans and monitoring turbidity levels.'
$ws.Cells.Item(11, 7).Value = 'Monitor and report on mitigation measures and submit reports by specified dates.'
$ws.Cells.Item(11, 8).Value = 'Offset habitat loss by placing course rock substrate and adhere to contingency plans.'
$ws.Cells.Item(11, 9).Value = 'Intertidal marine habitat, Subtidal marine habitat, Intertidal unnamed stream'
$ws.Cells.Item(11, 10).Value = 'None explicitly listed in the provided sections'
$ws.Cells.Item(11, 11).Value = '2792 HEUs of potential fish habitat'
$ws.Cells.Item(11, 12).Value = 'None mentioned'
$ws.Cells.Item(11, 13).Value = 'Course rock substrate placement as part of offset measures'
$ws.Cells.Item(11, 14).Value = 'None mentioned'
$ws.Cells.Item(11, 15).Value = 'None mentioned'
$ws.Cells.Item(11, 16).Value = $null

